$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Script1")

# Update credentials in the login data sheet
$ws.Range("A2").Value = "bhanu"
$ws.Range("B2").Value = "bhanu123"

# Update the active selection to match the edited workbook (B2 only)
$ws.Range("B2").Select()
